# Updated cryptos list - applies latest price/volume/name/link changes
# to the coin table (columns B=Coin, C=Link, D=Price, E=Volume(1h)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, [string]$text) {
    # Preserve the existing cell style while forcing the value to be
    # written as literal text, so numeric-looking strings (e.g. "19.74")
    # are not coerced into floating point numbers by Excel.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

# Row 2
Set-CellText $ws.Range("D2") '26.729.59'
Set-CellText $ws.Range("E2") '  +0.32%  '

# Row 3
Set-CellText $ws.Range("D3") '1.601.35'
Set-CellText $ws.Range("E3") '  +0.18%  '

# Row 4
Set-CellText $ws.Range("E4") '  +0.23%  '

# Row 5
Set-CellText $ws.Range("D5") '211.69'
Set-CellText $ws.Range("E5") '  +0.05%  '

# Row 6
Set-CellText $ws.Range("E6") '  -0.52%  '

# Row 7
Set-CellText $ws.Range("E7") '  +0.23%  '

# Row 8
Set-CellText $ws.Range("E8") '  +0.16%  '

# Row 9
Set-CellText $ws.Range("E9") '  +0.35%  '

# Row 10
Set-CellText $ws.Range("D10") '19.74'
Set-CellText $ws.Range("E10") '  +0.85%  '

# Row 11
Set-CellText $ws.Range("D11") '0.0846'
Set-CellText $ws.Range("E11") '  +0.87%  '

# Row 12
Set-CellText $ws.Range("D12") '1.826.06'
Set-CellText $ws.Range("E12") '  +0.17%  '

# Row 13
Set-CellText $ws.Range("D13") '1.622.21'
Set-CellText $ws.Range("E13") '  +1.44%  '

# Row 14
Set-CellText $ws.Range("E14") '  +0.47%  '

# Row 15
Set-CellText $ws.Range("E15") '  -0.09%  '

# Row 16
Set-CellText $ws.Range("D16") '65.04'
Set-CellText $ws.Range("E16") '  +0.06%  '

# Row 17
Set-CellText $ws.Range("B17") 'ShibaInu'
Set-CellText $ws.Range("C17") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText $ws.Range("D17") '0.0₃0740'
Set-CellText $ws.Range("E17") '  +0.83%  '

# Row 18
Set-CellText $ws.Range("B18") 'BitcoinCash'
Set-CellText $ws.Range("C18") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-CellText $ws.Range("D18") '210.02'
Set-CellText $ws.Range("E18") '  +0.91%  '

# Row 19
Set-CellText $ws.Range("B19") 'Chainlink'
Set-CellText $ws.Range("C19") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText $ws.Range("D19") '7.21'
Set-CellText $ws.Range("E19") '  +2.00%  '

# Row 20
Set-CellText $ws.Range("B20") 'Dai'
Set-CellText $ws.Range("C20") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText $ws.Range("D20") '1.01'
Set-CellText $ws.Range("E20") '  +0.28%  '

# Row 21
Set-CellText $ws.Range("B21") 'Uniswap'
Set-CellText $ws.Range("C21") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText $ws.Range("D21") '4.28'
Set-CellText $ws.Range("E21") '  -0.11%  '

# Row 22
Set-CellText $ws.Range("B22") 'Toncoin'
Set-CellText $ws.Range("C22") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-CellText $ws.Range("D22") '2.27'
Set-CellText $ws.Range("E22") '  -2.23%  '

# Row 23
Set-CellText $ws.Range("B23") 'Avalanche'
Set-CellText $ws.Range("C23") 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-CellText $ws.Range("D23") '8.97'
Set-CellText $ws.Range("E23") '  +0.35%  '

# Row 24
Set-CellText $ws.Range("B24") 'Monero'
Set-CellText $ws.Range("C24") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText $ws.Range("D24") '143.95'
Set-CellText $ws.Range("E24") '  -0.82%  '

# Row 25
Set-CellText $ws.Range("B25") 'BinanceUSD'
Set-CellText $ws.Range("C25") 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-CellText $ws.Range("D25") '1.01'
Set-CellText $ws.Range("E25") '  +0.13%  '

# Row 26
Set-CellText $ws.Range("B26") 'Cosmos'
Set-CellText $ws.Range("C26") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText $ws.Range("D26") '7.08'
Set-CellText $ws.Range("E26") '  -0.58%  '

# Row 27
Set-CellText $ws.Range("B27") 'Stellar'
Set-CellText $ws.Range("C27") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText $ws.Range("D27") '0.114'
Set-CellText $ws.Range("E27") '  -0.75%  '

# Row 28
Set-CellText $ws.Range("B28") 'EthereumClassic'
Set-CellText $ws.Range("C28") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-CellText $ws.Range("D28") '15.39'
Set-CellText $ws.Range("E28") '  +0.51%  '

# Row 29
Set-CellText $ws.Range("B29") 'Hedera'
Set-CellText $ws.Range("C29") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText $ws.Range("D29") '0.0511'
Set-CellText $ws.Range("E29") '  -0.30%  '

# Row 30
Set-CellText $ws.Range("B30") 'PancakeSwap'
Set-CellText $ws.Range("C30") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText $ws.Range("D30") '1.16'
Set-CellText $ws.Range("E30") '  +0.16%  '

# Row 31
Set-CellText $ws.Range("B31") 'Filecoin'
Set-CellText $ws.Range("C31") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws.Range("D31") '3.26'
Set-CellText $ws.Range("E31") '  +0.96%  '

# Row 32
Set-CellText $ws.Range("B32") 'InternetComputer(DFINITY)'
Set-CellText $ws.Range("C32") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws.Range("D32") '2.98'
Set-CellText $ws.Range("E32") '  +1.13%  '

# Row 33
Set-CellText $ws.Range("B33") 'Maker'
Set-CellText $ws.Range("C33") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText $ws.Range("D33") '1.294.82'
Set-CellText $ws.Range("E33") '  +1.33%  '

# Row 34
Set-CellText $ws.Range("B34") 'HuobiToken'
Set-CellText $ws.Range("C34") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText $ws.Range("D34") '2.47'
Set-CellText $ws.Range("E34") '  +0.77%  '

# Row 35
Set-CellText $ws.Range("B35") 'LidoDAOToken'
Set-CellText $ws.Range("C35") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-CellText $ws.Range("D35") '1.49'
Set-CellText $ws.Range("E35") '  +0.61%  '

# Row 36
Set-CellText $ws.Range("B36") 'ImmutableX'
Set-CellText $ws.Range("C36") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws.Range("D36") '0.596'
Set-CellText $ws.Range("E36") '  -4.06%  '

# Row 37
Set-CellText $ws.Range("B37") 'WEMIXToken'
Set-CellText $ws.Range("C37") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-CellText $ws.Range("D37") '1.16'
Set-CellText $ws.Range("E37") '  +11.73%  '

# Row 38
Set-CellText $ws.Range("B38") 'VeChain'
Set-CellText $ws.Range("C38") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws.Range("D38") '0.0169'
Set-CellText $ws.Range("E38") '  -0.78%  '

# Row 39
Set-CellText $ws.Range("B39") 'ARBITRUM'
Set-CellText $ws.Range("C39") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-CellText $ws.Range("D39") '0.829'
Set-CellText $ws.Range("E39") '  -0.73%  '

# Row 40
Set-CellText $ws.Range("B40") 'FraxShare'
Set-CellText $ws.Range("C40") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText $ws.Range("D40") '5.39'
Set-CellText $ws.Range("E40") '  -2.11%  '

# Row 41
Set-CellText $ws.Range("B41") 'MXToken'
Set-CellText $ws.Range("C41") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws.Range("D41") '2.20'
Set-CellText $ws.Range("E41") '  -0.09%  '

# Row 42
Set-CellText $ws.Range("B42") 'TrustWalletToken'
Set-CellText $ws.Range("C42") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText $ws.Range("D42") '0.780'
Set-CellText $ws.Range("E42") '  -0.53%  '

# Row 43
Set-CellText $ws.Range("B43") 'Aave'
Set-CellText $ws.Range("C43") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText $ws.Range("D43") '63.01'
Set-CellText $ws.Range("E43") '  -1.49%  '

# Row 44
Set-CellText $ws.Range("B44") 'RocketPoolETH'
Set-CellText $ws.Range("C44") 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-CellText $ws.Range("D44") '1.738.35'
Set-CellText $ws.Range("E44") '  +0.24%  '

# Row 45
Set-CellText $ws.Range("B45") 'Quant'
Set-CellText $ws.Range("C45") 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText $ws.Range("D45") '90.59'
Set-CellText $ws.Range("E45") '  -0.09%  '

# Row 46
Set-CellText $ws.Range("B46") 'RenderToken'
Set-CellText $ws.Range("C46") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText $ws.Range("D46") '1.56'
Set-CellText $ws.Range("E46") '  -2.50%  '

# Row 47
Set-CellText $ws.Range("B47") 'Algorand'
Set-CellText $ws.Range("C47") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText $ws.Range("D47") '0.101'
Set-CellText $ws.Range("E47") '  -0.11%  '

# Row 48
Set-CellText $ws.Range("B48") 'Cronos'
Set-CellText $ws.Range("C48") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws.Range("D48") '0.0517'
Set-CellText $ws.Range("E48") '  +1.68%  '

# Row 49
Set-CellText $ws.Range("B49") 'USDD'
Set-CellText $ws.Range("C49") 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-CellText $ws.Range("D49") '1.01'
Set-CellText $ws.Range("E49") '  +0.31%  '

# Row 50
Set-CellText $ws.Range("B50") 'EnergySwap'
Set-CellText $ws.Range("C50") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText $ws.Range("D50") '7.41'
Set-CellText $ws.Range("E50") '  +0.23%  '

# Row 51
Set-CellText $ws.Range("B51") 'Mantle'
Set-CellText $ws.Range("C51") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-CellText $ws.Range("D51") '0.398'
Set-CellText $ws.Range("E51") '  +1.27%  '
